# "Generate Report for Handoff" — refresh the handoff timestamps for the
# most-recently-handed-off file (f8b920fb-b355-4bba-b9ea-2e9686a1e0b9)
# across the Overview summary sheet and the per-locale detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) for the
#     f8b920fb... row (row 5) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-09-09 11:54:55"

# --- zh-cn detail sheet: "Latest Handoff Datetime" (column H) for the
#     f8b920fb... row (row 5) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-09-09 11:54:44"

# --- de-de detail sheet: "Latest Handoff Datetime" (column H) for the
#     f8b920fb... row (row 5) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-09-09 11:54:55"
